$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")

# Insert a new "unit" column before the existing "equipmentClass" column:
# copy the old D1 header ("equipmentClass") into the new E1 cell, then
# overwrite D1 with the new header "unit".
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("D1").Value = "unit"

# Fill the new "unit" column (D2:D67) with the SI unit identifiers used
# throughout the DCC schema (base units, derived units, accepted non-SI
# units, and prefixes).
$units = @(
    "\metre",
    "\kilogram",
    "\second",
    "\ampere",
    "\kelvin",
    "\mole",
    "\candela",
    "\one",
    "\day",
    "\hour",
    "\minute",
    "\degree",
    "\arcminute",
    "\arcsecond",
    "\gram",
    "\radian",
    "\steradian",
    "\hertz",
    "\newton",
    "\pascal",
    "\joule",
    "\watt",
    "\coulomb",
    "\volt",
    "\farad",
    "\ohm",
    "\siemens",
    "\weber",
    "\tesla",
    "\henry",
    "\degreecelsius",
    "\lumen",
    "\lux",
    "\becquerel",
    "\sievert",
    "\gray",
    "\katal",
    "\hectare",
    "\litre",
    "\tonne",
    "\electronvolt",
    "\dalton",
    "\astronomicalunit",
    "\neper",
    "\bel",
    "\decibel",
    "\deca",
    "\hecto",
    "\kilo",
    "\mega",
    "\giga",
    "\tera",
    "\peta",
    "\exa",
    "\zetta",
    "\yotta",
    "\deci",
    "\centi",
    "\milli",
    "\micro",
    "\nano",
    "\pico",
    "\femto",
    "\atto",
    "\zepto",
    "\yocto"
)

for ($i = 0; $i -lt $units.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $units[$i]
}

# Restore the active selection to where editing left off.
[void]$ws.Range("E4").Select()

$ws2 = $wb.Worksheets.Item("Table2")
[void]$ws2.Range("A6").Select()
